# Auto-generated data refresh for Siren_Profits workbook
# Applies updated currentAveragePrice* / Leve* derived columns (H:N) per sheet/row
$wb = $excel.ActiveWorkbook

# ALC!row80
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80,8).Value = 61257
$ws.Cells.Item(80,9).Value = 91620.77
$ws.Cells.Item(80,10).Value = 529.4545000000001
$ws.Cells.Item(80,11).Value = 274862.31
$ws.Cells.Item(80,12).Value = 1588.3635
$ws.Cells.Item(80,13).Value = -273864.31
$ws.Cells.Item(80,14).Value = -3584.3635

# ALC!row83
$ws.Cells.Item(83,8).Value = 61257
$ws.Cells.Item(83,9).Value = 91620.77
$ws.Cells.Item(83,10).Value = 529.4545000000001
$ws.Cells.Item(83,11).Value = 824586.9300000001
$ws.Cells.Item(83,12).Value = 4765.0905
$ws.Cells.Item(83,13).Value = -819594.9300000001
$ws.Cells.Item(83,14).Value = -14749.0905

# ALC!row92
$ws.Cells.Item(92,8).Value = 600
$ws.Cells.Item(92,9).Value = 644.17645
$ws.Cells.Item(92,11).Value = 644.17645
$ws.Cells.Item(92,13).Value = 603.82355

# ALC!row113
$ws.Cells.Item(113,8).Value = 18585.857
$ws.Cells.Item(113,10).Value = 15233
$ws.Cells.Item(113,12).Value = 15233
$ws.Cells.Item(113,14).Value = -21741

# ALC!row132
$ws.Cells.Item(132,8).Value = 2781370.2
$ws.Cells.Item(132,9).Value = 3509.5881
$ws.Cells.Item(132,11).Value = 10528.7643
$ws.Cells.Item(132,13).Value = -7998.764299999999

# ARM!row28
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28,8).Value = 173013330
$ws.Cells.Item(28,9).Value = 19986.5
$ws.Cells.Item(28,11).Value = 19986.5
$ws.Cells.Item(28,13).Value = -19794.5

# ARM!row32
$ws.Cells.Item(32,8).Value = 4660.6665
$ws.Cells.Item(32,9).Value = 4342.836
$ws.Cells.Item(32,10).Value = 8919.6
$ws.Cells.Item(32,11).Value = 4342.836
$ws.Cells.Item(32,12).Value = 8919.6
$ws.Cells.Item(32,13).Value = -4055.836
$ws.Cells.Item(32,14).Value = -9493.6

# ARM!row45
$ws.Cells.Item(45,8).Value = 148156
$ws.Cells.Item(45,9).Value = 226467.11
$ws.Cells.Item(45,10).Value = 7196
$ws.Cells.Item(45,11).Value = 226467.11
$ws.Cells.Item(45,12).Value = 7196
$ws.Cells.Item(45,13).Value = -226090.11
$ws.Cells.Item(45,14).Value = -7950

# ARM!row74
$ws.Cells.Item(74,8).Value = 7060
$ws.Cells.Item(74,10).Value = 2518.182
$ws.Cells.Item(74,12).Value = 2518.182
$ws.Cells.Item(74,14).Value = -4266.182

# ARM!row77
$ws.Cells.Item(77,8).Value = 7060
$ws.Cells.Item(77,10).Value = 2518.182
$ws.Cells.Item(77,12).Value = 12590.91
$ws.Cells.Item(77,14).Value = -21326.91

# ARM!row80
$ws.Cells.Item(80,8).Value = 77071.28999999999
$ws.Cells.Item(80,10).Value = 84916.5
$ws.Cells.Item(80,12).Value = 84916.5
$ws.Cells.Item(80,14).Value = -86912.5

# ARM!row83
$ws.Cells.Item(83,8).Value = 77071.28999999999
$ws.Cells.Item(83,10).Value = 84916.5
$ws.Cells.Item(83,12).Value = 254749.5
$ws.Cells.Item(83,14).Value = -264733.5

# ARM!row99
$ws.Cells.Item(99,8).Value = 173013330
$ws.Cells.Item(99,9).Value = 19986.5
$ws.Cells.Item(99,11).Value = 19986.5
$ws.Cells.Item(99,13).Value = -16991.5

# ARM!row122
$ws.Cells.Item(122,8).Value = 1161042.1
$ws.Cells.Item(122,9).Value = 6631.0625
$ws.Cells.Item(122,11).Value = 19893.1875
$ws.Cells.Item(122,13).Value = -17443.1875

# ARM!row132
$ws.Cells.Item(132,8).Value = 2357.2593
$ws.Cells.Item(132,9).Value = 1444.9143
$ws.Cells.Item(132,11).Value = 4334.742899999999
$ws.Cells.Item(132,13).Value = -1804.742899999999

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107,8).Value = 2236.5
$ws.Cells.Item(107,9).Value = 2528.2727
$ws.Cells.Item(107,11).Value = 2528.2727
$ws.Cells.Item(107,13).Value = -608.2727

# CRP!row9
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9,8).Value = 16500
$ws.Cells.Item(9,10).Value = 16500
$ws.Cells.Item(9,12).Value = 16500
$ws.Cells.Item(9,14).Value = -16836

# CRP!row16
$ws.Cells.Item(16,8).Value = 732.625
$ws.Cells.Item(16,9).Value = 714.9231
$ws.Cells.Item(16,11).Value = 714.9231
$ws.Cells.Item(16,13).Value = -427.9231

# CRP!row31
$ws.Cells.Item(31,8).Value = 6634.154
$ws.Cells.Item(31,9).Value = 6976.727
$ws.Cells.Item(31,10).Value = 4750
$ws.Cells.Item(31,11).Value = 6976.727
$ws.Cells.Item(31,12).Value = 4750
$ws.Cells.Item(31,13).Value = -6681.727
$ws.Cells.Item(31,14).Value = -5340

# CRP!row34
$ws.Cells.Item(34,8).Value = 6634.154
$ws.Cells.Item(34,9).Value = 6976.727
$ws.Cells.Item(34,10).Value = 4750
$ws.Cells.Item(34,11).Value = 6976.727
$ws.Cells.Item(34,12).Value = 4750
$ws.Cells.Item(34,13).Value = -6774.727
$ws.Cells.Item(34,14).Value = -5154

# CRP!row68
$ws.Cells.Item(68,8).Value = 0
$ws.Cells.Item(68,10).Value = 0
$ws.Cells.Item(68,12).Value = 0
$ws.Cells.Item(68,14).ClearContents()

# CRP!row71
$ws.Cells.Item(71,8).Value = 0
$ws.Cells.Item(71,10).Value = 0
$ws.Cells.Item(71,12).Value = 0
$ws.Cells.Item(71,14).ClearContents()

# CRP!row107
$ws.Cells.Item(107,8).Value = 10022.143
$ws.Cells.Item(107,9).Value = 14345.889
$ws.Cells.Item(107,10).Value = 2239.4
$ws.Cells.Item(107,11).Value = 14345.889
$ws.Cells.Item(107,12).Value = 2239.4
$ws.Cells.Item(107,13).Value = -12425.889
$ws.Cells.Item(107,14).Value = -6079.4

# CRP!row113
$ws.Cells.Item(113,8).Value = 732.625
$ws.Cells.Item(113,9).Value = 714.9231
$ws.Cells.Item(113,11).Value = 714.9231
$ws.Cells.Item(113,13).Value = 1455.0769

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12,8).Value = 48.944443
$ws.Cells.Item(12,9).Value = 59.3
$ws.Cells.Item(12,10).Value = 36
$ws.Cells.Item(12,11).Value = 177.9
$ws.Cells.Item(12,12).Value = 108
$ws.Cells.Item(12,13).Value = -4.899999999999977
$ws.Cells.Item(12,14).Value = -454

# CUL!row33
$ws.Cells.Item(33,8).Value = 53.714287
$ws.Cells.Item(33,10).Value = 55.25
$ws.Cells.Item(33,12).Value = 331.5
$ws.Cells.Item(33,14).Value = -897.5

# CUL!row49
$ws.Cells.Item(49,8).Value = 0
$ws.Cells.Item(49,9).Value = 0
$ws.Cells.Item(49,10).Value = 0
$ws.Cells.Item(49,11).Value = 0
$ws.Cells.Item(49,12).Value = 0
$ws.Cells.Item(49,13).ClearContents()
$ws.Cells.Item(49,14).ClearContents()

# CUL!row51
$ws.Cells.Item(51,8).Value = 1050

# CUL!row55
$ws.Cells.Item(55,8).Value = 5017.4
$ws.Cells.Item(55,10).Value = 5979.4165
$ws.Cells.Item(55,12).Value = 17938.2495
$ws.Cells.Item(55,14).Value = -18292.2495

# CUL!row94
$ws.Cells.Item(94,8).Value = 900
$ws.Cells.Item(94,9).Value = 900
$ws.Cells.Item(94,11).Value = 2700
$ws.Cells.Item(94,13).Value = -2024

# CUL!row121
$ws.Cells.Item(121,8).Value = 3178.2856
$ws.Cells.Item(121,9).Value = 3350
$ws.Cells.Item(121,10).Value = 3049.5
$ws.Cells.Item(121,11).Value = 10050
$ws.Cells.Item(121,12).Value = 9148.5
$ws.Cells.Item(121,13).Value = -8740
$ws.Cells.Item(121,14).Value = -11768.5

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102,8).Value = 18372.555
$ws.Cells.Item(102,9).Value = 19919.125
$ws.Cells.Item(102,11).Value = 19919.125
$ws.Cells.Item(102,13).Value = -18297.125

# GSM!row122
$ws.Cells.Item(122,8).Value = 10953.546
$ws.Cells.Item(122,9).Value = 6769.5884
$ws.Cells.Item(122,10).Value = 25179
$ws.Cells.Item(122,11).Value = 20308.7652
$ws.Cells.Item(122,12).Value = 75537
$ws.Cells.Item(122,13).Value = -17858.7652
$ws.Cells.Item(122,14).Value = -80437

# GSM!row132
$ws.Cells.Item(132,8).Value = 6219.7393
$ws.Cells.Item(132,9).Value = 6592.421
$ws.Cells.Item(132,10).Value = 4449.5
$ws.Cells.Item(132,11).Value = 19777.263
$ws.Cells.Item(132,12).Value = 13348.5
$ws.Cells.Item(132,13).Value = -17247.263
$ws.Cells.Item(132,14).Value = -18408.5

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55,8).Value = 1977.4445
$ws.Cells.Item(55,9).Value = 314.5
$ws.Cells.Item(55,11).Value = 314.5
$ws.Cells.Item(55,13).Value = -141.5

# LTW!row82
$ws.Cells.Item(82,8).Value = 2412.625
$ws.Cells.Item(82,9).Value = 3268.3333
$ws.Cells.Item(82,10).Value = 1312.4286
$ws.Cells.Item(82,11).Value = 3268.3333
$ws.Cells.Item(82,12).Value = 1312.4286
$ws.Cells.Item(82,13).Value = -2907.3333
$ws.Cells.Item(82,14).Value = -2034.4286

# LTW!row85
$ws.Cells.Item(85,8).Value = 2412.625
$ws.Cells.Item(85,9).Value = 3268.3333
$ws.Cells.Item(85,10).Value = 1312.4286
$ws.Cells.Item(85,11).Value = 3268.3333
$ws.Cells.Item(85,12).Value = 1312.4286
$ws.Cells.Item(85,13).Value = -2020.3333
$ws.Cells.Item(85,14).Value = -3808.4286

# LTW!row93
$ws.Cells.Item(93,8).Value = 7034.5713
$ws.Cells.Item(93,9).Value = 8136.8184
$ws.Cells.Item(93,10).Value = 2993
$ws.Cells.Item(93,11).Value = 8136.8184
$ws.Cells.Item(93,12).Value = 2993
$ws.Cells.Item(93,13).Value = -6888.8184
$ws.Cells.Item(93,14).Value = -5489

# LTW!row132
$ws.Cells.Item(132,8).Value = 442113.06
$ws.Cells.Item(132,9).Value = 787345.3
$ws.Cells.Item(132,10).Value = 4818.8667
$ws.Cells.Item(132,11).Value = 2362035.9
$ws.Cells.Item(132,12).Value = 14456.6001
$ws.Cells.Item(132,13).Value = -2359505.9
$ws.Cells.Item(132,14).Value = -19516.6001

# LTW!row136
$ws.Cells.Item(136,8).Value = 6074.8335

# WVR!row70
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70,8).Value = 482137.6
$ws.Cells.Item(70,10).Value = 771266.3
$ws.Cells.Item(70,12).Value = 771266.3
$ws.Cells.Item(70,14).Value = -771896.3

# WVR!row73
$ws.Cells.Item(73,8).Value = 482137.6
$ws.Cells.Item(73,10).Value = 771266.3
$ws.Cells.Item(73,12).Value = 771266.3
$ws.Cells.Item(73,14).Value = -773450.3

# WVR!row113
$ws.Cells.Item(113,8).Value = 1890.9565
$ws.Cells.Item(113,9).Value = 985.5
$ws.Cells.Item(113,11).Value = 2956.5
$ws.Cells.Item(113,13).Value = -786.5

# WVR!row126
$ws.Cells.Item(126,8).Value = 15916.182
$ws.Cells.Item(126,9).Value = 20520.959
$ws.Cells.Item(126,10).Value = 3636.7778
$ws.Cells.Item(126,11).Value = 61562.87699999999
$ws.Cells.Item(126,12).Value = 10910.3334
$ws.Cells.Item(126,13).Value = -59092.87699999999
$ws.Cells.Item(126,14).Value = -15850.3334

# WVR!row132
$ws.Cells.Item(132,8).Value = 3414.25
$ws.Cells.Item(132,9).Value = 2758.4211
$ws.Cells.Item(132,10).Value = 4798.778
$ws.Cells.Item(132,11).Value = 8275.263300000001
$ws.Cells.Item(132,12).Value = 14396.334
$ws.Cells.Item(132,13).Value = -5745.263300000001
$ws.Cells.Item(132,14).Value = -19456.334

# WVR!row136
$ws.Cells.Item(136,8).Value = 1735622.5
$ws.Cells.Item(136,9).Value = 1940075.2
$ws.Cells.Item(136,11).Value = 5820225.6
$ws.Cells.Item(136,13).Value = -5817675.6

